# "first revision of demo mode"
# The Air Pressure sensor byte readout changed from a float encoding to an
# int encoding, so relabel the four "Air Pressure Float Byte N" rows on the
# "Bytes" sheet to "Air Pressure Int Byte N".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bytes")

$ws.Range("C12").Value = "Air Pressure Int Byte 0"
$ws.Range("C13").Value = "Air Pressure Int Byte 1"
$ws.Range("C14").Value = "Air Pressure Int Byte 2"
$ws.Range("C15").Value = "Air Pressure Int Byte 3"

# Restore the view: scroll back to the top of the sheet and move the
# selection from D25 up to A5.
$ws.Range("A1").Select()
$ws.Range("A5").Select()

# Best-effort resize of the saved window geometry to match the author's
# larger editing window.
try { $wb.Windows.Item(1).Width = 51200 } catch {}
try { $wb.Windows.Item(1).Height = 28260 } catch {}
try { $excel.ActiveWindow.Width = 51200 } catch {}
try { $excel.ActiveWindow.Height = 28260 } catch {}
